$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: STM32C011 -> STM32C031, price/qty update, link update ---
$ws.Range("B7").Value = "STM32C031"
$ws.Range("D7").Value = 6.67
$ws.Range("E7").Formula = "=D7*C7"
$ws.Range("F7").Value = "https://pl.rs-online.com/web/p/mikrokontrolery/0214865?gb=s"

# Re-typing B7's text resets its "quote prefix" cell style (s=2 -> s=1);
# restore it by pulling the format back from A7, which carries the same style.
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# --- F4: drop the "price from contest rules - " prefix, keep just the URL ---
# (use .Formula rather than .Value here so the existing Hyperlink cell style, s=8, survives)
$ws.Range("F4").Formula = "https://www.digikey.pl/pl/products/detail/texas-instruments/TLV2462CDGKR/1677686"

# --- Clear the "price from contest rules" placeholder text in F2, F3, F5, F8 ---
$ws.Range("F2").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("F8").ClearContents()

# --- Rebuild hyperlinks: drop the stale display text / duplicate relationship
#     by removing the old hyperlink first, then re-adding clean ones so rIds
#     line up as rId1=digikey(F4), rId2=rs-online(F7), rId3=botland(F6) ---
$ws.Range("F4").Hyperlinks.Delete()
[void]$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.digikey.pl/pl/products/detail/texas-instruments/TLV2462CDGKR/1677686")
[void]$ws.Hyperlinks.Add($ws.Range("F7"), "https://pl.rs-online.com/web/p/mikrokontrolery/0214865?gb=s")
[void]$ws.Hyperlinks.Add($ws.Range("F6"), "https://botland.com.pl/akcesoria-do-raspberry-pi-pico/18854-zestaw-zlacz-meskich-do-gpio-raspberry-pi-pico-5904422328511.html")

# --- Give F6 and F7 the same "Hyperlink" look (underline + theme color + border) F4 already has ---
$ws.Range("F4").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the now-obsolete blank rows above the total, shifting the spacer
#     row (13) and the total row (14) up to rows 10 and 11 ---
[void]$ws.Range("A9:A11").EntireRow.Delete()

# --- View-state cosmetics to mirror the saved workbook ---
[void]$ws.Range("B17").Select()

Write-Output "done"
